# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The workbook lists workers with outstanding balances ("Estado de Cuenta").
# The refresh removes the stale trailing periods for ANGEL MANUEL ROMERO COTA
# (2201-2208, all sharing one "Valor Mora") and replaces them with the
# updated periods/values, while MARCO JOSE PADILLA MELENDEZ now carries the
# most recent period (2208) and ANDERSON DE JESUS CABRALES PADILLA is
# re-appended at the bottom of the table with his original period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: MARCO JOSE PADILLA MELENDEZ -> period 2208, new date/valor mora
$ws.Range("C19").Value = "73147453"
$ws.Range("D19").Value = "MARCO JOSE PADILLA MELENDEZ"
$ws.Range("E19").Value = "2208"
$ws.Range("F19").Value = 25749
$ws.Range("G19").Value = 908526

# Rows 20-26: ANGEL MANUEL ROMERO COTA, periods 2207 down to 2201, new valor mora
$ws.Range("C20").Value = "73208312"
$ws.Range("D20").Value = "ANGEL MANUEL ROMERO COTA"
$ws.Range("E20").Value = "2207"
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 908526

$ws.Range("C21").Value = "73208312"
$ws.Range("D21").Value = "ANGEL MANUEL ROMERO COTA"
$ws.Range("E21").Value = "2206"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 908526

$ws.Range("C22").Value = "73208312"
$ws.Range("D22").Value = "ANGEL MANUEL ROMERO COTA"
$ws.Range("E22").Value = "2205"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 908526

$ws.Range("C23").Value = "73208312"
$ws.Range("D23").Value = "ANGEL MANUEL ROMERO COTA"
$ws.Range("E23").Value = "2204"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 908526

$ws.Range("C24").Value = "73208312"
$ws.Range("D24").Value = "ANGEL MANUEL ROMERO COTA"
$ws.Range("E24").Value = "2203"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 908526

$ws.Range("C25").Value = "73208312"
$ws.Range("D25").Value = "ANGEL MANUEL ROMERO COTA"
$ws.Range("E25").Value = "2202"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 908526

$ws.Range("C26").Value = "73208312"
$ws.Range("D26").Value = "ANGEL MANUEL ROMERO COTA"
$ws.Range("E26").Value = "2201"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 908526

# Row 27: ANDERSON DE JESUS CABRALES PADILLA re-added at the bottom, period 2006
$ws.Range("C27").Value = "1047436296"
$ws.Range("D27").Value = "ANDERSON DE JESUS CABRALES PADILLA"
$ws.Range("E27").Value = "2006"
$ws.Range("F27").Value = 35112
$ws.Range("G27").Value = 877803

# Re-fit the data columns to the refreshed (longer) values, matching the
# author's re-save of the sheet.
$ws.Columns.Item(2).ColumnWidth = 17.6
$ws.Columns.Item(3).ColumnWidth = 15.77
$ws.Columns.Item(5).ColumnWidth = 12.6
$ws.Columns.Item(6).ColumnWidth = 9.26
$ws.Columns.Item(7).ColumnWidth = 13.42
$ws.Columns.Item(8).ColumnWidth = 18.42
$ws.Columns.Item(9).ColumnWidth = 17.27
$ws.Columns.Item(10).ColumnWidth = 14.1
